# Auto-generated Excel COM-interop script
# Applies the scheduled market-data refresh captured in the commit diff to the
# Leve-profit workbook (columns H:N -- price / profit columns -- across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117.413795
$ws.Range("I33").Value = 73.708336
$ws.Range("K33").Value = 73.708336
$ws.Range("M33").Value = 155.291664
$ws.Range("H40").Value = 3491.2727
$ws.Range("I40").Value = 3875.25
$ws.Range("J40").Value = 2467.3333
$ws.Range("K40").Value = 3875.25
$ws.Range("L40").Value = 2467.3333
$ws.Range("M40").Value = -3700.25
$ws.Range("N40").Value = -2817.3333
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H53").Value = 321.66666
$ws.Range("I53").Value = 198.8
$ws.Range("J53").Value = 567.4
$ws.Range("K53").Value = 198.8
$ws.Range("L53").Value = 567.4
$ws.Range("M53").Value = 438.2
$ws.Range("N53").Value = -1841.4
$ws.Range("H70").Value = 2546.6
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2546.6
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = 7639.799999999999
$ws.Range("N70").Value = -8179.799999999999
$ws.Range("H73").Value = 2546.6
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2546.6
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = 7639.799999999999
$ws.Range("N73").Value = -9511.799999999999
$ws.Range("H132").Value = 212503
$ws.Range("I132").Value = 243587.3
$ws.Range("J132").Value = 39812.445
$ws.Range("K132").Value = 730761.8999999999
$ws.Range("L132").Value = 119437.335
$ws.Range("M132").Value = -728231.8999999999
$ws.Range("N132").Value = -124497.335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3043.4626
$ws.Range("I32").Value = 1819.125
$ws.Range("K32").Value = 1819.125
$ws.Range("M32").Value = -1532.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 49
$ws.Range("I22").Value = 49
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 49
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 15025000
$ws.Range("I6").Value = 15025000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 15025000
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -15024887
$ws.Range("H31").Value = 1326.7179
$ws.Range("I31").Value = 915.13043
$ws.Range("J31").Value = 1918.375
$ws.Range("K31").Value = 915.13043
$ws.Range("L31").Value = 1918.375
$ws.Range("M31").Value = -620.13043
$ws.Range("N31").Value = -2508.375
$ws.Range("H34").Value = 1326.7179
$ws.Range("I34").Value = 915.13043
$ws.Range("J34").Value = 1918.375
$ws.Range("K34").Value = 915.13043
$ws.Range("L34").Value = 1918.375
$ws.Range("M34").Value = -713.13043
$ws.Range("N34").Value = -2322.375
$ws.Range("H58").Value = 1985.6842
$ws.Range("I58").Value = 624.75
$ws.Range("J58").Value = 2975.4546
$ws.Range("K58").Value = 624.75
$ws.Range("L58").Value = 2975.4546
$ws.Range("M58").Value = -421.75
$ws.Range("N58").Value = -3381.4546
$ws.Range("H99").Value = 6945842
$ws.Range("I99").Value = 10417927
$ws.Range("K99").Value = 10417927
$ws.Range("M99").Value = -10416429
$ws.Range("H115").Value = 26782.25
$ws.Range("J115").Value = 26782.25
$ws.Range("L115").Value = 26782.25
$ws.Range("N115").Value = -29132.25
$ws.Range("H120").Value = 33662.668
$ws.Range("J120").Value = 33662.668
$ws.Range("L120").Value = 33662.668
$ws.Range("N120").Value = -40920.668
$ws.Range("H126").Value = 6945842
$ws.Range("I126").Value = 10417927
$ws.Range("K126").Value = 31253781
$ws.Range("M126").Value = -31251311
$ws.Range("H136").Value = 1985.6842
$ws.Range("I136").Value = 624.75
$ws.Range("J136").Value = 2975.4546
$ws.Range("K136").Value = 1874.25
$ws.Range("L136").Value = 8926.363799999999
$ws.Range("M136").Value = 675.75
$ws.Range("N136").Value = -14026.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 815.22
$ws.Range("I68").Value = 692.058
$ws.Range("J68").Value = 1089.3549
$ws.Range("K68").Value = 2076.174
$ws.Range("L68").Value = 3268.0647
$ws.Range("M68").Value = -1265.174
$ws.Range("N68").Value = -4890.0647
$ws.Range("H71").Value = 815.22
$ws.Range("I71").Value = 692.058
$ws.Range("J71").Value = 1089.3549
$ws.Range("K71").Value = 6228.522
$ws.Range("L71").Value = 9804.194100000001
$ws.Range("M71").Value = -2172.522
$ws.Range("N71").Value = -17916.1941
$ws.Range("H106").Value = 2750
$ws.Range("I106").Value = 2500
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 7500
$ws.Range("L106").Value = 9000
$ws.Range("M106").Value = -6554
$ws.Range("N106").Value = -10892
$ws.Range("H107").Value = 739.8182
$ws.Range("I107").Value = 164.73914
$ws.Range("J107").Value = 1153.1562
$ws.Range("K107").Value = 494.2174199999999
$ws.Range("L107").Value = 3459.4686
$ws.Range("M107").Value = 1425.78258
$ws.Range("N107").Value = -7299.4686
$ws.Range("H129").Value = 712
$ws.Range("I129").Value = 440
$ws.Range("J129").Value = 1800
$ws.Range("K129").Value = 1320
$ws.Range("L129").Value = 5400
$ws.Range("M129").Value = 3680
$ws.Range("N129").Value = -15400
$ws.Range("H131").Value = 2429.0952
$ws.Range("I131").Value = 516.44446
$ws.Range("J131").Value = 2658.6133
$ws.Range("K131").Value = 1549.33338
$ws.Range("L131").Value = 7975.8399
$ws.Range("M131").Value = 3490.66662
$ws.Range("N131").Value = -18055.8399
$ws.Range("H134").Value = 19086.912
$ws.Range("I134").Value = 2219.9
$ws.Range("J134").Value = 32061.54
$ws.Range("K134").Value = 6659.700000000001
$ws.Range("L134").Value = 96184.62
$ws.Range("M134").Value = -1589.700000000001
$ws.Range("N134").Value = -106324.62
$ws.Range("H137").Value = 4212997.5
$ws.Range("I137").Value = 6670594
$ws.Range("J137").Value = 117003.664
$ws.Range("K137").Value = 20011782
$ws.Range("L137").Value = 351010.992
$ws.Range("M137").Value = -20006682
$ws.Range("N137").Value = -361210.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1853435.1
$ws.Range("I122").Value = 3704703.8
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 11114111.4
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -11111661.4
$ws.Range("N122").Value = -11400.0001
$ws.Range("H126").Value = 2280.75
$ws.Range("I126").Value = 1903
$ws.Range("J126").Value = 2469.625
$ws.Range("K126").Value = 5709
$ws.Range("L126").Value = 7408.875
$ws.Range("M126").Value = -3239
$ws.Range("N126").Value = -12348.875
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3291.6667
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 3545.4546
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 3545.4546
$ws.Range("M7").Value = -388
$ws.Range("N7").Value = -3769.4546
$ws.Range("H14").Value = 12196.2
$ws.Range("I14").Value = 26000.5
$ws.Range("J14").Value = 2993.3333
$ws.Range("K14").Value = 26000.5
$ws.Range("L14").Value = 2993.3333
$ws.Range("M14").Value = -25828.5
$ws.Range("N14").Value = -3337.3333
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H63").Value = 32000
$ws.Range("J63").Value = 32000
$ws.Range("L63").Value = 32000
$ws.Range("N63").Value = -33498
$ws.Range("H66").Value = 32000
$ws.Range("J66").Value = 32000
$ws.Range("L66").Value = 96000
$ws.Range("N66").Value = -103488
$ws.Range("H126").Value = 3291.6667
$ws.Range("I126").Value = 500
$ws.Range("J126").Value = 3545.4546
$ws.Range("K126").Value = 1500
$ws.Range("L126").Value = 10636.3638
$ws.Range("M126").Value = 970
$ws.Range("N126").Value = -15576.3638
$ws.Range("H136").Value = 3546.389
$ws.Range("I136").Value = 1138.72
$ws.Range("J136").Value = 9018.362999999999
$ws.Range("K136").Value = 3416.16
$ws.Range("L136").Value = 27055.089
$ws.Range("M136").Value = -866.1599999999999
$ws.Range("N136").Value = -32155.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 11118914
$ws.Range("J24").Value = 8778.375
$ws.Range("L24").Value = 8778.375
$ws.Range("N24").Value = -9238.375
$ws.Range("H26").Value = 43333.332
$ws.Range("H122").Value = 335001.34
$ws.Range("I122").Value = 335001.34
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1005004.02
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1002554.02
$ws.Range("H126").Value = 100779.9
$ws.Range("I126").Value = 111755.445
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 335266.335
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -332796.335
$ws.Range("N126").Value = -10940
$ws.Range("H136").Value = 20897522
$ws.Range("I136").Value = 30394536
$ws.Range("J136").Value = 4090.6
$ws.Range("K136").Value = 91183608
$ws.Range("L136").Value = 12271.8
$ws.Range("M136").Value = -91181058
$ws.Range("N136").Value = -17371.8

